$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Cell D2 held "Name_AlreadyExist_Search" - update it to the new search term used by the test.
# Keep the trailing spaces - the saved workbook stores them with xml:space="preserve".
$ws.Range("D2").Value = "selenium_GF1   "

# The selection moved from C2 to D11 when the workbook was last saved.
$ws.Range("D11").Select()
